$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the header style already used by G1 ("sum")
# by copying G1's formatting so the new cell reuses the same style record.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
